# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" right after "2021-Q4" (before "总计"),
#    using a copy of "2021-Q4" as a template so it inherits the same
#    column layout / header styling.
# 2. Fill "2022-Q1" with the new fund-holding snapshot.
# 3. Prepend a "2022-Q1" row to the "总计" (totals) summary sheet, pushing
#    the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet as a copy of "2021-Q4", positioned
# immediately after it (i.e. before "总计").
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($null, $templateSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------
# Step 2: overwrite the copied data with the 2022-Q1 fund snapshot.
# Header row (B1:H1) is already correct from the template copy, just
# like the rest of the sheet formatting / styles.
# Numeric-looking text fields are written with a leading apostrophe so
# they stay text (matching the existing sheet convention) instead of
# being auto-converted to numbers.
# ---------------------------------------------------------------------
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'164811"
$newSheet.Range("C2").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）A"
$newSheet.Range("D2").Value = "'0.23"
$newSheet.Range("E2").Value = "'94.28"
$newSheet.Range("F2").Value = "'3.56"
$newSheet.Range("G2").Value = "'0.0082"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'512780"
$newSheet.Range("C3").Value = "广发中证京津冀协同发展主题ETF"
$newSheet.Range("D3").Value = "'0.13"
$newSheet.Range("E3").Value = "'98.52"
$newSheet.Range("F3").Value = "'3.18"
$newSheet.Range("G3").Value = "'0.0041"
$newSheet.Range("H3").Value = 5

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'164825"
$newSheet.Range("C4").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）C"
$newSheet.Range("D4").Value = "'0.06"
$newSheet.Range("E4").Value = "'94.28"
$newSheet.Range("F4").Value = "'3.56"
$newSheet.Range("G4").Value = "'0.0021"
$newSheet.Range("H4").Value = 4

# ---------------------------------------------------------------------
# Step 3: update the "总计" (totals) sheet — push the two existing data
# rows down by one, then write the new first row for "2022-Q1".
# (Values are written literally rather than copied cell-to-cell, since
# the COM Range.Value *getter* is unreliable in this host - Range.Value2
# works for reads, but the final values here are known up front anyway.)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Row 4 is brand new, so it has no style yet - copy column A's styling
# (bold + border) from row 2 so it matches rows 2/3.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

# Row 4 <- old row 3 (2021-Q2).
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q2"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 0.01

# Row 3 <- old row 2 (2021-Q4).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.01

# New first data row: 2022-Q1.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.01

# ---------------------------------------------------------------------
# Restore the original active sheet/tab - copying a sheet makes the new
# copy the active one, but "2021-Q2" was the active tab originally.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
